$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (preserve exact formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "45.968.98"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "2.333.76"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "299.37"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").Value = "98.14"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.507"
$ws.Range("E9").Value = "  -6.02%  "
$ws.Range("D10").Value = "34.39"
$ws.Range("E10").Value = "  -4.48%  "
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -3.32%  "
$ws.Range("D12").Value = "7.07"
$ws.Range("E12").Value = "  -5.39%  "
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "2.689.42"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "2.335.13"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "13.58"
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("D17").Value = "0.800"
$ws.Range("E17").Value = "  -4.64%  "
$ws.Range("D18").Value = "45.898.46"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  -8.57%  "
$ws.Range("D20").Value = "0.0₃0959"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").Value = "5.93"
$ws.Range("E21").Value = "  -4.75%  "
$ws.Range("D22").Value = "66.19"
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").Value = "242.70"
$ws.Range("E23").Value = "  -3.57%  "
$ws.Range("D24").Value = "2.79"
$ws.Range("E24").Value = "  -6.72%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -6.01%  "
$ws.Range("D27").Value = "40.23"
$ws.Range("E27").Value = "  -5.43%  "
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("D29").Value = "9.59"
$ws.Range("E29").Value = "  -3.55%  "
$ws.Range("D30").Value = "20.61"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "3.59"
$ws.Range("E31").Value = "  +14.03%  "
$ws.Range("E32").Value = "  +6.37%  "
$ws.Range("D33").Value = "5.39"
$ws.Range("E33").Value = "  -7.85%  "
$ws.Range("D34").Value = "144.06"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("D35").Value = "0.0767"
$ws.Range("E35").Value = "  -6.20%  "
$ws.Range("D36").Value = "0.111"
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("E37").Value = "  -3.75%  "
$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").Value = "15.16"
$ws.Range("E39").Value = "  +8.07%  "
$ws.Range("D40").Value = "3.83"
$ws.Range("E40").Value = "  -4.41%  "
$ws.Range("D41").Value = "0.0296"
$ws.Range("E41").Value = "  -5.09%  "
$ws.Range("D42").Value = "3.15"
$ws.Range("E42").Value = "  -7.44%  "
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "1.847.49"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("D45").Value = "90.38"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").Value = "1.79"
$ws.Range("E46").Value = "  -9.46%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.184"
$ws.Range("E47").Value = "  -6.34%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "70.17"
$ws.Range("E48").Value = "  -6.30%  "
$ws.Range("D49").Value = "2.561.85"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "95.52"
$ws.Range("E50").Value = "  -3.59%  "
$ws.Range("D51").Value = "4.71"
$ws.Range("E51").Value = "  -1.53%  "
